# Update cryptos list data (price + 1h volume change) per upstream refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the literal text into the cell (even when it "looks" numeric,
    # e.g. "592.34" or "1.00") without leaving a formula behind and without
    # introducing a new cell style (which a plain .Value assignment of a
    # numeric-looking string would do, since Excel auto-converts it to a
    # Number). We build a literal-text formula, then copy/paste-special as
    # values so the result is a plain inline/shared text cell, matching the
    # original workbooks cell representation.
    $escaped = $text.Replace("""", """""")
    $cell.Formula = "=""" + $escaped + """"
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$ws.Application.CutCopyMode = $false

Set-TextValue $ws.Cells.Item(2, 4) "61.182.69"
$ws.Cells.Item(2, 5).Value = "  +0.94%  "
Set-TextValue $ws.Cells.Item(3, 4) "2.931.45"
$ws.Cells.Item(3, 5).Value = "  +0.96%  "
$ws.Cells.Item(4, 5).Value = "  +0.06%  "
Set-TextValue $ws.Cells.Item(5, 4) "592.34"
$ws.Cells.Item(5, 5).Value = "  +0.87%  "
Set-TextValue $ws.Cells.Item(6, 4) "146.31"
$ws.Cells.Item(6, 5).Value = "  -1.13%  "
$ws.Cells.Item(7, 5).Value = "  +0.04%  "
$ws.Cells.Item(8, 5).Value = "  +0.32%  "
Set-TextValue $ws.Cells.Item(9, 4) "6.88"
$ws.Cells.Item(9, 5).Value = "  +2.64%  "
$ws.Cells.Item(10, 5).Value = "  +0.09%  "
$ws.Cells.Item(11, 5).Value = "  -1.26%  "
$ws.Cells.Item(12, 5).Value = "  +1.10%  "
Set-TextValue $ws.Cells.Item(13, 4) "33.85"
$ws.Cells.Item(13, 5).Value = "  -1.27%  "
$ws.Cells.Item(14, 5).Value = "  -0.62%  "
Set-TextValue $ws.Cells.Item(15, 4) "3.417.40"
$ws.Cells.Item(15, 5).Value = "  +1.06%  "
Set-TextValue $ws.Cells.Item(16, 4) "61.196.50"
$ws.Cells.Item(16, 5).Value = "  +1.05%  "
$ws.Cells.Item(17, 5).Value = "  -1.62%  "
Set-TextValue $ws.Cells.Item(18, 4) "2.926.26"
$ws.Cells.Item(18, 5).Value = "  +0.92%  "
Set-TextValue $ws.Cells.Item(19, 4) "432.33"
$ws.Cells.Item(19, 5).Value = "  +1.50%  "
Set-TextValue $ws.Cells.Item(20, 4) "13.46"
$ws.Cells.Item(20, 5).Value = "  -1.43%  "
Set-TextValue $ws.Cells.Item(21, 4) "0.684"
$ws.Cells.Item(21, 5).Value = "  +1.72%  "
Set-TextValue $ws.Cells.Item(22, 4) "7.10"
$ws.Cells.Item(22, 5).Value = "  -0.40%  "
Set-TextValue $ws.Cells.Item(23, 4) "81.43"
$ws.Cells.Item(23, 5).Value = "  +0.87%  "
Set-TextValue $ws.Cells.Item(24, 4) "11.06"
$ws.Cells.Item(24, 5).Value = "  -0.23%  "
$ws.Cells.Item(25, 5).Value = "  +1.17%  "
$ws.Cells.Item(26, 5).Value = "  +1.81%  "
$ws.Cells.Item(27, 5).Value = "  -0.04%  "
Set-TextValue $ws.Cells.Item(28, 4) "2.32"
$ws.Cells.Item(28, 5).Value = "  +5.80%  "
Set-TextValue $ws.Cells.Item(29, 4) "1.00"
$ws.Cells.Item(29, 5).Value = "  +0.04%  "
$ws.Cells.Item(30, 5).Value = "  +0.07%  "
Set-TextValue $ws.Cells.Item(31, 4) "7.11"
$ws.Cells.Item(31, 5).Value = "  -2.67%  "
Set-TextValue $ws.Cells.Item(32, 4) "26.52"
$ws.Cells.Item(32, 5).Value = "  -0.10%  "
Set-TextValue $ws.Cells.Item(33, 4) "0.108"
$ws.Cells.Item(33, 5).Value = "  +1.57%  "
$ws.Cells.Item(34, 5).Value = "  +3.10%  "
$ws.Cells.Item(35, 5).Value = "  +0.58%  "
Set-TextValue $ws.Cells.Item(36, 4) "5.63"
$ws.Cells.Item(36, 5).Value = "  -0.82%  "
Set-TextValue $ws.Cells.Item(37, 4) "3.10"
$ws.Cells.Item(37, 5).Value = "  +3.79%  "
Set-TextValue $ws.Cells.Item(38, 4) "50.00"
$ws.Cells.Item(38, 5).Value = "  +1.34%  "
$ws.Cells.Item(39, 5).Value = "  +3.20%  "
$ws.Cells.Item(40, 5).Value = "  -1.45%  "
Set-TextValue $ws.Cells.Item(41, 4) "8.59"
$ws.Cells.Item(41, 5).Value = "  -1.82%  "
$ws.Cells.Item(42, 5).Value = "  -0.94%  "
$ws.Cells.Item(43, 5).Value = "  -4.69%  "
Set-TextValue $ws.Cells.Item(44, 4) "376.38"
$ws.Cells.Item(44, 5).Value = "  +1.25%  "
$ws.Cells.Item(45, 5).Value = "  +0.24%  "
Set-TextValue $ws.Cells.Item(46, 4) "2.715.61"
$ws.Cells.Item(46, 5).Value = "  +2.33%  "
Set-TextValue $ws.Cells.Item(47, 4) "130.75"
$ws.Cells.Item(47, 5).Value = "  -1.84%  "
$ws.Cells.Item(48, 5).Value = "  +0.02%  "
Set-TextValue $ws.Cells.Item(49, 4) "24.30"
$ws.Cells.Item(49, 5).Value = "  -3.10%  "
$ws.Cells.Item(50, 5).Value = "  +0.19%  "
$ws.Cells.Item(51, 5).Value = "  -2.70%  "

$ws.Application.CutCopyMode = $false
